# Daily auto push: insert a new reading (2026/02/04, 23:00, rank 23 -> 201)
# into the time series sheet. The new row is inserted above row 782, which
# shifts all the existing rows 782:823 down to 783:824 (matching how the
# sheet grows by one data point per push).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above row 782; everything below (old rows 782:823)
# shifts down to 783:824, and the sheet's dimension grows to A1:D824
# automatically.
$ws.Rows("782:782").Insert()

# Column A holds a date stored as text (e.g. "2026/02/04"). Typing that
# string straight into a General-formatted cell would make Excel reinterpret
# it as a real date serial number, which is not what the source data uses.
# Copying it from an existing cell that already holds the exact same text
# (row 780 is also 2026/02/04) keeps it as plain text with the sheet's
# default (unstyled) formatting, just like the rest of the column.
$ws.Range("A780").Copy($ws.Range("A782"))

$ws.Range("B782").Value = "水"
$ws.Range("C782").Value = 23
$ws.Range("D782").Value = 201
